$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($cellRef).Style = "Normal"
}

Set-TextCell "D2" "35.043.44"
Set-TextCell "E2" "  +1.14%  "

Set-TextCell "D3" "1.847.12"
Set-TextCell "E3" "  +2.14%  "

Set-TextCell "E4" "  +0.01%  "

Set-TextCell "D5" "233.28"
Set-TextCell "E5" "  +0.70%  "

Set-TextCell "D6" "0.622"
Set-TextCell "E6" "  +3.09%  "

Set-TextCell "E7" "  -0.08%  "

Set-TextCell "E8" "  +7.25%  "

Set-TextCell "D9" "0.328"
Set-TextCell "E9" "  +0.28%  "

Set-TextCell "E10" "  +2.06%  "

Set-TextCell "E11" "  -0.80%  "

Set-TextCell "D12" "2.115.32"
Set-TextCell "E12" "  +2.11%  "

Set-TextCell "D13" "11.35"
Set-TextCell "E13" "  +1.79%  "

Set-TextCell "B14" "Polygon"
Set-TextCell "C14" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextCell "D14" "0.677"
Set-TextCell "E14" "  +0.47%  "

Set-TextCell "B15" "WrappedEther"
Set-TextCell "C15" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell "D15" "1.841.23"
Set-TextCell "E15" "  +1.78%  "

Set-TextCell "D16" "4.69"
Set-TextCell "E16" "  +2.38%  "

Set-TextCell "D17" "35.045.24"
Set-TextCell "E17" "  +1.16%  "

Set-TextCell "D18" "69.97"
Set-TextCell "E18" "  +0.82%  "

Set-TextCell "D19" "0.0₃0793"
Set-TextCell "E19" "  +1.20%  "

Set-TextCell "D20" "240.81"
Set-TextCell "E20" "  +0.21%  "

Set-TextCell "D21" "12.11"
Set-TextCell "E21" "  +1.91%  "

Set-TextCell "E22" "  +1.98%  "

Set-TextCell "E23" "  -0.09%  "

Set-TextCell "E24" "  +2.86%  "

Set-TextCell "D25" "171.68"
Set-TextCell "E25" "  +0.11%  "

Set-TextCell "D26" "7.88"
Set-TextCell "E26" "  +2.27%  "

Set-TextCell "B27" "PancakeSwap"
Set-TextCell "C27" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell "D27" "1.77"
Set-TextCell "E27" "  +15.40%  "

Set-TextCell "B28" "EthereumClassic"
Set-TextCell "C28" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell "D28" "17.60"
Set-TextCell "E28" "  +2.39%  "

Set-TextCell "B29" "Stellar"
Set-TextCell "C29" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell "D29" "0.124"
Set-TextCell "E29" "  +3.29%  "

Set-TextCell "B30" "BinanceUSD"
Set-TextCell "C30" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextCell "D30" "1.01"
Set-TextCell "E30" "  +0.04%  "

Set-TextCell "B31" "Hedera"
Set-TextCell "C31" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D31" "0.0556"
Set-TextCell "E31" "  +2.09%  "

Set-TextCell "E32" "  -1.34%  "

Set-TextCell "E33" "  +0.83%  "

Set-TextCell "E34" "  +24.43%  "

Set-TextCell "D35" "1.98"
Set-TextCell "E35" "  +10.78%  "

Set-TextCell "D36" "0.765"
Set-TextCell "E36" "  +9.43%  "

Set-TextCell "D37" "1.23"
Set-TextCell "E37" "  -5.68%  "

Set-TextCell "D38" "1.07"
Set-TextCell "E38" "  +11.71%  "

Set-TextCell "D39" "91.57"
Set-TextCell "E39" "  +0.16%  "

Set-TextCell "D40" "0.0201"
Set-TextCell "E40" "  +4.66%  "

Set-TextCell "D41" "1.348.53"
Set-TextCell "E41" "  +1.64%  "

Set-TextCell "D42" "14.66"
Set-TextCell "E42" "  +3.53%  "

Set-TextCell "E43" "  +5.14%  "

Set-TextCell "D44" "12.69"
Set-TextCell "E44" "  +85.42%  "

Set-TextCell "E45" "  -3.29%  "

Set-TextCell "E46" "  +3.22%  "

Set-TextCell "D47" "0.0531"
Set-TextCell "E47" "  +3.64%  "

Set-TextCell "E48" "  +1.09%  "

Set-TextCell "D49" "2.028.40"
Set-TextCell "E49" "  +1.48%  "

Set-TextCell "E50" "  +17.34%  "

Set-TextCell "E51" "  +0.69%  "

